$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New weekly crime-data collection: refresh the period headers plus every
# count / percentage cell in the Murder..Traffic Fatalities table (rows
# 14-33) with this week's figures.
# ---------------------------------------------------------------------------

# Header text: "Volume 32   Number  24" -> "...25" and the covered week's
# date range moves forward by one week.
$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"

# Bulk numeric updates for the statistics grid (Week to Date / 28 Day /
# Year to Date / 2 Year columns and their % Chg figures).
$data = @(
  @("D14", 1),
  @("E14", 200),
  @("F14", 8),
  @("G14", 8),
  @("H14", 0),
  @("I14", 47),
  @("J14", 53),
  @("K14", -11.320754716981),
  @("L14", -26.5625),
  @("M14", -12.962962962963),
  @("N14", -79.017857142857),
  @("C15", 10),
  @("D15", 13),
  @("E15", -23.076923076923),
  @("F15", 43),
  @("G15", 42),
  @("H15", 2.380952380952),
  @("I15", 249),
  @("J15", 209),
  @("K15", 19.138755980861),
  @("L15", 33.155080213903),
  @("M15", 85.820895522388),
  @("N15", -26.331360946745),
  @("C16", 105),
  @("D16", 104),
  @("E16", 0.961538461538),
  @("F16", 402),
  @("G16", 382),
  @("H16", 5.235602094240),
  @("I16", 2174),
  @("J16", 2275),
  @("K16", -4.439560439560),
  @("L16", 2.547169811320),
  @("M16", 9.908998988877),
  @("N16", -71.224354731965),
  @("C17", 158),
  @("D17", 206),
  @("E17", -23.300970873786),
  @("F17", 743),
  @("G17", 745),
  @("H17", -0.268456375838),
  @("I17", 3987),
  @("J17", 3814),
  @("K17", 4.535920293654),
  @("L17", 9.292763157894),
  @("M17", 95.441176470588),
  @("N17", -3.181155900922),
  @("C18", 63),
  @("D18", 55),
  @("E18", 14.545454545454),
  @("F18", 208),
  @("G18", 215),
  @("H18", -3.255813953488),
  @("I18", 1341),
  @("J18", 1358),
  @("K18", -1.251840942562),
  @("L18", -6.875),
  @("M18", -9.757738896366),
  @("N18", -84.768287142208),
  @("C19", 169),
  @("D19", 165),
  @("E19", 2.424242424242),
  @("F19", 702),
  @("G19", 690),
  @("H19", 1.739130434782),
  @("I19", 4201),
  @("J19", 4252),
  @("K19", -1.199435559736),
  @("L19", 14.562312517043),
  @("M19", 103.044949250846),
  @("N19", 23.449897149573),
  @("C20", 75),
  @("D20", 85),
  @("E20", -11.764705882352),
  @("F20", 331),
  @("G20", 328),
  @("H20", 0.914634146341),
  @("I20", 2087),
  @("J20", 1903),
  @("K20", 9.668943772990),
  @("L20", -17.215390717969),
  @("M20", 120.380147835269),
  @("N20", -71.185972663261),
  @("C21", 583),
  @("D21", 629),
  @("E21", -7.313195548489),
  @("F21", 2437),
  @("G21", 2410),
  @("H21", 1.120331950207),
  @("I21", 14086),
  @("J21", 13864),
  @("K21", 1.601269474899),
  @("L21", 3.216824210449),
  @("M21", 61.759301791456),
  @("N21", -55.543632633738),
  @("C22", 4),
  @("D22", 2),
  @("E22", 100),
  @("F22", 25),
  @("G22", 17),
  @("H22", 47.058823529411),
  @("I22", 143),
  @("J22", 159),
  @("K22", -10.062893081761),
  @("L22", -4.026845637583),
  @("M22", -11.728395061728),
  @("C23", 36),
  @("D23", 25),
  @("E23", 44),
  @("F23", 128),
  @("G23", 121),
  @("H23", 5.785123966942),
  @("I23", 730),
  @("J23", 789),
  @("K23", -7.477820025348),
  @("L23", -12.259615384615),
  @("M23", 53.361344537815),
  @("C24", 335),
  @("D24", 281),
  @("E24", 19.217081850533),
  @("F24", 1455),
  @("G24", 1080),
  @("H24", 34.722222222222),
  @("I24", 8364),
  @("J24", 7614),
  @("K24", 9.850275807722),
  @("L24", -0.381133873272),
  @("M24", 45.157931273863),
  @("C25", 111),
  @("D25", 97),
  @("E25", 14.432989690721),
  @("F25", 481),
  @("G25", 407),
  @("H25", 18.181818181818),
  @("I25", 2739),
  @("J25", 3033),
  @("K25", -9.693372898120),
  @("L25", -24.938339271033),
  @("C26", 234),
  @("D26", 253),
  @("E26", -7.509881422924),
  @("F26", 963),
  @("G26", 988),
  @("H26", -2.530364372469),
  @("I26", 5181),
  @("J26", 5192),
  @("K26", -0.211864406779),
  @("L26", 4.645526156332),
  @("M26", 0.038617493724),
  @("C27", 13),
  @("D27", 20),
  @("E27", -35),
  @("F27", 55),
  @("G27", 58),
  @("H27", -5.172413793103),
  @("I27", 315),
  @("J27", 323),
  @("K27", -2.476780185758),
  @("L27", 0),
  @("C28", 16),
  @("D28", 27),
  @("E28", -40.740740740740),
  @("F28", 97),
  @("G28", 106),
  @("H28", -8.490566037735),
  @("I28", 521),
  @("J28", 580),
  @("K28", -10.172413793103),
  @("L28", 2.964426877470),
  @("C29", 4),
  @("D29", 15),
  @("E29", -73.333333333333),
  @("G29", 37),
  @("H29", -29.729729729729),
  @("I29", 135),
  @("J29", 176),
  @("K29", -23.295454545454),
  @("L29", -18.181818181818),
  @("M29", -31.472081218274),
  @("N29", -76.764199655765),
  @("C30", 4),
  @("D30", 12),
  @("E30", -66.666666666666),
  @("G30", 31),
  @("H30", -25.806451612903),
  @("I30", 117),
  @("J30", 143),
  @("K30", -18.181818181818),
  @("L30", -16.428571428571),
  @("M30", -30.357142857142),
  @("N30", -77.756653992395),
  @("L31", -18.181818181818),
  @("F33", 1),
  @("H33", -83.333333333333),
  @("L33", -42.857142857142)
)

foreach ($pair in $data) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# C33 (Traffic Fatalities, Week-to-Date 2025 column) has no fatalities this
# week, so the cell switches from the numeric value 1 to the text "0" used
# elsewhere in the sheet for a suppressed/non-applicable count - matching the
# style already used by its text-holding neighbor D33.
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "0"
$ws.Range("D33").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$excel.CutCopyMode = $false
